# test5.xlsx — replace the "Nombre/Apellidos/.../pollingStation" person
# record with the smaller "Nombre/locacalizacion/Correo electrónico/id/kind"
# location record (ExcelParseTest / TextParseTest fixture update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "locacalizacion"
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("D1").Value = "id"
$ws.Range("E1").Value = "kind"

# --- Data row ---------------------------------------------------------------
$ws.Range("A2").Value = "jorge"
$ws.Range("B2").Value = "18:13:14:12S"
$ws.Range("C2").Value = "jorge@email.es"
$ws.Range("D2").Value = "ID4"
$ws.Range("E2").Value = 1

# The old sheet used columns A:I — the new layout only needs A:E, so drop
# the now-unused trailing columns (Nacionalidad, DNI, NIF, pollingStation).
$ws.Range("F1:I2").ClearContents()

# The old C2 (email) cell carried a real mailto: hyperlink; the new data
# keeps the hyperlink-style formatting on C2 but no longer has a live link.
foreach ($hl in $ws.Hyperlinks) {
    $hl.Delete()
}

# Match the saved selection/used-range of the authored workbook.
$ws.Range("A1:E2").Select()
